$excel | Get-Member | Out-String | Write-Output
